$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.022883666666667
$ws.Range("H2").Value = 3.068651
$ws.Range("I2").Value = 0.1594660351460709
$ws.Range("J2").Value = 0.1594660351460709
$ws.Range("M2").Value = 1.485259333333333
$ws.Range("N2").Value = 4.455778
$ws.Range("O2").Value = 0.3057455162066235
$ws.Range("P2").Value = 0.3057455162066235
$ws.Range("Q2").Value = 1.519247512830889
$ws.Range("R2").Value = 13.673227615478
$ws.Range("S2").Value = 0.04875602523315902
$ws.Range("T2").Value = 0.04875602523315901
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.022883666666667
$ws.Range("H3").Value = 3.068651
$ws.Range("I3").Value = 0.1594660351460709
$ws.Range("J3").Value = 0.1594660351460709
$ws.Range("O3").Value = 0.2805555239151429
$ws.Range("P3").Value = 0.2805555239151429
$ws.Range("Q3").Value = 1.394078602385778
$ws.Range("R3").Value = 12.546707421472
$ws.Range("S3").Value = 0.04473907703707653
$ws.Range("T3").Value = 0.04473907703707652
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.022883666666667
$ws.Range("H4").Value = 3.068651
$ws.Range("I4").Value = 0.1594660351460709
$ws.Range("J4").Value = 0.1594660351460709
$ws.Range("O4").Value = 0.4136989598782336
$ws.Range("P4").Value = 0.4136989598782336
$ws.Range("Q4").Value = 2.055667483381778
$ws.Range("R4").Value = 18.501007350436
$ws.Range("S4").Value = 0.0659709328758354
$ws.Range("T4").Value = 0.06597093287583539
$ws.Range("I5").Value = 0.1603506552336246
$ws.Range("J5").Value = 0.1603506552336246
$ws.Range("M5").Value = 1.485259333333333
$ws.Range("N5").Value = 4.455778
$ws.Range("O5").Value = 0.3057455162066235
$ws.Range("P5").Value = 0.3057455162066235
$ws.Range("Q5").Value = 1.527675369374667
$ws.Range("R5").Value = 13.749078324372
$ws.Range("S5").Value = 0.04902649385847486
$ws.Range("T5").Value = 0.04902649385847486
$ws.Range("I6").Value = 0.1603506552336246
$ws.Range("J6").Value = 0.1603506552336246
$ws.Range("O6").Value = 0.2805555239151429
$ws.Range("P6").Value = 0.2805555239151429
$ws.Range("S6").Value = 0.044987262089206
$ws.Range("T6").Value = 0.044987262089206
$ws.Range("I7").Value = 0.1603506552336246
$ws.Range("J7").Value = 0.1603506552336246
$ws.Range("O7").Value = 0.4136989598782336
$ws.Range("P7").Value = 0.4136989598782336
$ws.Range("S7").Value = 0.06633689928594373
$ws.Range("T7").Value = 0.06633689928594373
$ws.Range("I8").Value = 0.6801833096203046
$ws.Range("J8").Value = 0.6801833096203045
$ws.Range("M8").Value = 1.485259333333333
$ws.Range("N8").Value = 4.455778
$ws.Range("O8").Value = 0.3057455162066235
$ws.Range("P8").Value = 0.3057455162066235
$ws.Range("Q8").Value = 6.480168648221335
$ws.Range("R8").Value = 58.32151783399201
$ws.Range("S8").Value = 0.2079629971149896
$ws.Range("T8").Value = 0.2079629971149896
$ws.Range("I9").Value = 0.6801833096203046
$ws.Range("J9").Value = 0.6801833096203045
$ws.Range("O9").Value = 0.2805555239151429
$ws.Range("P9").Value = 0.2805555239151429
$ws.Range("S9").Value = 0.1908291847888604
$ws.Range("T9").Value = 0.1908291847888604
$ws.Range("I10").Value = 0.6801833096203046
$ws.Range("J10").Value = 0.6801833096203045
$ws.Range("O10").Value = 0.4136989598782336
$ws.Range("P10").Value = 0.4136989598782336
$ws.Range("S10").Value = 0.2813911277164546
$ws.Range("T10").Value = 0.2813911277164545